$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header cells for Wins/Losses/Ties, matching the existing header formatting
$ws.Range("AC1").Copy()
$ws.Range("AD1:AF1").PasteSpecial(-4122)
$ws.Range("AD1").Value = "Wins"
$ws.Range("AE1").Value = "Losses"
$ws.Range("AF1").Value = "Ties"

# Team record data for every data row (2-57)
for ($r = 2; $r -le 57; $r++) {
    $ws.Cells.Item($r, 30).Value = 83  # AD - Wins
    $ws.Cells.Item($r, 31).Value = 79  # AE - Losses
    $ws.Cells.Item($r, 32).Value = 0   # AF - Ties
}
